$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 207, shifting existing rows 207:282 down to 208:283.
$ws.Rows("207:207").Insert()

# Populate the newly inserted row 207 with its data.
$ws.Range("A207").Value = 5
$ws.Range("B207").Value = "Macroferia Regional de Talca"
$ws.Range("C207").Value = "Maule"
$ws.Range("D207").Value = 44468
$ws.Range("E207").Value = 7
$ws.Range("F207").Value = 100114001
$ws.Range("G207").Value = "Papa"
$ws.Range("H207").Value = "Rodeo"
$ws.Range("I207").Value = "1a (guarda lavada)"
$ws.Range("J207").Value = 1500
$ws.Range("K207").Value = 10000
$ws.Range("L207").Value = 10000
$ws.Range("M207").Value = 10000
$ws.Range("N207").Value = "$/malla 25 kilos"
$ws.Range("O207").Value = "Región de Los Lagos"
$ws.Range("P207").Value = 400
$ws.Range("Q207").Value = 25
$ws.Range("R207").Value = "Hortaliza"
